# 13_LibFormula.xlsx - "configurazione br e moduli nuovi 18-5-2017"
#
# Add a new indicator row (INDICATOR_119) to the Library_Formula sheet,
# inserted right after the INDICATOR_118 row (row 105), pushing the
# subsequent rows (INDICATOR_150, INDICATOR_151, INDICATOR_104,
# INDICATOR_900) down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Library_Formula")

# Insert a brand-new row at 105 (shifts existing rows 105-108 down to 106-109)
$ws.Rows.Item(105).Insert() | Out-Null

# Populate the new row the same way every other indicator row in this
# table is populated: Action / Library / Formula Name / Output / Input
$ws.Range("A105").Value = "CREATE/MODIFY"
$ws.Range("B105").Value = "LIB_EWS_RETAIL"
$ws.Range("C105").Value = "INDICATOR_119"
$ws.Range("E105").Value = "String"
$ws.Range("F105").Value = "String"

# Leave the selection where the author left it after the edit
$ws.Range("A104:F105").Select() | Out-Null
